# Apply the "global-transaction -> statistics" menu consolidation edit.
#
# IAM_PERMISSION: the site-level "global-transaction" saga / saga-instance
#   route rows (11:12) are removed outright - the equivalent "statistics"
#   routes already exist a couple of rows below and shift up to take
#   their place.
#
# IAM_MENU_B: the site-level "global-transaction" menu group and its two
#   children (rows 12:14) are removed outright - the equivalent
#   "statistics" menu group (with its own children) already exists below
#   and shifts up to take their place.
#
# IAM_MENU_PERMISSION: the "global-transaction.saga" / "global-transaction
#   .saga-instance" menu-permission rows (17:22) are renamed in place to
#   "statistics.saga" / "statistics.saga-instance" (the PERMISSION_CODE
#   values in column G are unaffected), and the now-duplicate old
#   "statistics.saga" / "statistics.saga-instance" rows (23:28) are
#   removed outright.

$wb = $excel.ActiveWorkbook

# --- IAM_PERMISSION ---------------------------------------------------
$wsPermission = $wb.Worksheets.Item("IAM_PERMISSION")
$wsPermission.Rows("11:12").Delete()

# --- IAM_MENU_B ---------------------------------------------------------
$wsMenu = $wb.Worksheets.Item("IAM_MENU_B")
$wsMenu.Rows("12:14").Delete()

# --- IAM_MENU_PERMISSION -------------------------------------------------
$wsMenuPermission = $wb.Worksheets.Item("IAM_MENU_PERMISSION")
$wsMenuPermission.Range("F17").Value = "choerodon.code.statistics.saga"
$wsMenuPermission.Range("F18").Value = "choerodon.code.statistics.saga"
$wsMenuPermission.Range("F19").Value = "choerodon.code.statistics.saga-instance"
$wsMenuPermission.Range("F20").Value = "choerodon.code.statistics.saga-instance"
$wsMenuPermission.Range("F21").Value = "choerodon.code.statistics.saga-instance"
$wsMenuPermission.Range("F22").Value = "choerodon.code.statistics.saga-instance"
$wsMenuPermission.Rows("23:28").Delete()
